$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-02-16 Friday" "2024-02-17 Saturday"

Replace-Text "507×3=" "982×5="
Replace-Text "424×7=" "341×2="
Replace-Text "953×8=" "290×4="
Replace-Text "892×7=" "340×2="
Replace-Text "846×8=" "943×6="

Replace-Text "493×8=" "504×5="
Replace-Text "976×4=" "973×2="
Replace-Text "994×6=" "673×9="
Replace-Text "579×7=" "863×7="
Replace-Text "523×3=" "794×9="

Replace-Text "288×5=" "841×3="
Replace-Text "674×6=" "753×9="
Replace-Text "461×3=" "317×4="
Replace-Text "767×9=" "696×7="
Replace-Text "505×5=" "611×8="

Replace-Text "683×2=" "512×5="
Replace-Text "508×5=" "486×6="
Replace-Text "882×7=" "448×9="
Replace-Text "999×7=" "792×9="
Replace-Text "431×3=" "148×3="

Replace-Text "540×5=" "241×4="
Replace-Text "858×4=" "365×9="
Replace-Text "490×3=" "793×4="
Replace-Text "177×8=" "870×8="
Replace-Text "942×3=" "217×7="
